$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 458 in the
# "Feria Lagunitas de Puerto Montt - Cebollín" sheet, pushing the
# previously existing rows 458-534 down to 459-535.
$ws.Rows.Item(458).Insert()

$ws.Range("A458").Value = 4
$ws.Range("B458").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C458").Value = "Los Lagos"
$ws.Range("D458").Value = 45218
$ws.Range("E458").Value = 10
$ws.Range("F458").Value = 100112037
$ws.Range("G458").Value = "Cebollín"
$ws.Range("H458").Value = "Sin especificar"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 70
$ws.Range("K458").Value = 6500
$ws.Range("L458").Value = 6500
$ws.Range("M458").Value = 6500
$ws.Range("N458").Value = "$/paquete 36 unidades"
$ws.Range("O458").Value = "Región Metropolitana"
$ws.Range("P458").Value = 181
$ws.Range("Q458").Value = 36
$ws.Range("R458").Value = "Hortaliza"
